# Finished UDP Sender example
$wb = $excel.ActiveWorkbook

$wsJson = $wb.Worksheets.Item("JSON Reader Code Beispiel")
$wsUdp  = $wb.Worksheets.Item("UDP Sender Code Beispiel")

# --- JSON Reader sheet: add COUNTIF formulas (values stay the same: 10 / 15) ---
$wsJson.Range("B32").Formula = "=COUNTIF(D2:D26,""Yes"")"
$wsJson.Range("B33").Formula = "=COUNTIF(D2:D26,""No"")"

# --- JSON Reader sheet: selection moved from F34 to C36 ---
$wsJson.Range("C36").Select() | Out-Null

# --- UDP Sender sheet: fill in the CompilerErr (B) / Index (C) columns for rows 7-26 ---
$udpData = @{
    7  = @(61, 0)
    8  = @(63, 0)
    9  = @(69, 0)
    10 = @(56, 0)
    11 = @(61, 0)
    12 = @(58, 1)
    13 = @(55, 0)
    14 = @(65, 0)
    15 = @(60, 0)
    16 = @(67, 0)
    17 = @(56, 0)
    18 = @(65, 0)
    19 = @(65, 0)
    20 = @(75, 0)
    21 = @(66, 0)
    22 = @(60, 0)
    23 = @(56, 0)
    24 = @(63, 0)
    25 = @(65, 1)
    26 = @(56, 0)
}

foreach ($r in $udpData.Keys) {
    $vals = $udpData[$r]
    $wsUdp.Range("B$r").Value = $vals[0]
    $wsUdp.Range("C$r").Value = $vals[1]
}

# --- UDP Sender sheet: add COUNTIF formulas for Anzahl Yes / Anzahl No ---
$wsUdp.Range("B32").Formula = "=COUNTIF(D2:D26,""Yes"")"
$wsUdp.Range("B33").Formula = "=COUNTIF(D2:D26,""No"")"

# --- UDP Sender sheet: selection moved from F6 to B34 ---
$wsUdp.Range("B34").Select() | Out-Null

$excel.Calculate() | Out-Null
